$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.281367421150208
$ws.Range("B1").Value = 2.217790365219116
$ws.Range("C1").Value = 2.750782251358032
$ws.Range("D1").Value = 3.177966356277466
$ws.Range("E1").Value = 2.338335037231445
